$d = $word.ActiveDocument

$bodyInner = @'
<w:p w:rsidR="00003F14" w:rsidRDefault="007D48DD">
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="44"/>
        </w:rPr>
        <w:t>Weaving the Threads of Human Existence: Government, Society, and the Interplay of Power</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00003F14" w:rsidRDefault="007D48DD">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="36"/>
        </w:rPr>
        <w:t>Professor Alexander Graham</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00003F14" w:rsidRDefault="007D48DD">
      <w:pPr>
        <w:jc w:val="center"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>alexandergraham83@topmail</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
        </w:rPr>
        <w:t>edu</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00003F14" w:rsidRDefault="00003F14"/>
    <w:p w:rsidR="00003F14" w:rsidRDefault="007D48DD">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>In the intricate fabric of human society, nations and their governing systems stand as cornerstones, shaping the tapestry of our lives</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> From the ancient civilizations of yore to the complex political landscapes of today, understanding how governments function, how societies evolve, and how power dynamics interact is fundamental to our comprehension of the world around us</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Governments act as stewards of our collective interests, implementing policies that determine the well-being of their citizens</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> They are responsible for maintaining law and order, regulating economic activity, providing public services, and upholding the rights and freedoms of individuals</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> The structures, systems, and processes of governance vary across countries and cultures, yet the underlying principles of leadership, representation, and accountability remain the same</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Governments can be classified into different types, including democracies, monarchies, and autocracies</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Each type has its own unique characteristics and mechanisms for exercising power</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>The interplay of power is a pervasive force within societies, as individuals, groups, and institutions vie for influence and control</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Power can be wielded through various means, including economic resources, political authority, military might, cultural hegemony, and social networks</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> The struggle for power can lead to conflicts and tensions within societies but can also drive social and political change</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Understanding the dynamics of power is essential for analyzing political processes, social inequalities, and historical events</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> It also helps us appreciate the importance of checks and balances, transparency, and ethical leadership in preventing power from being abused</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Societies, as dynamic and evolving entities, are shaped by a multitude of factors, including culture, religion, ethnicity, language, and economic systems</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> They possess a unique identity and a shared sense of history, values, and norms</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> As societies evolve, they often face challenges such as demographic changes, technological transformations, environmental issues, and global interconnectedness</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Understanding the complexities of societies requires an </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>interdisciplinary approach that draws upon fields such as sociology, economics, anthropology, and history</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00003F14" w:rsidRDefault="007D48DD">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Summary</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00003F14" w:rsidRDefault="007D48DD">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>In this essay, we explored the fascinating world of government, society, and the intricate interplay of power</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> We examined the multifaceted roles of governments in shaping the lives of their citizens</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> We analyzed the different types of governments and the mechanisms through which they exercise power</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> Furthermore, we delved into the dynamics of power, its manifestations within societies, and the importance of checks and balances in ensuring ethical governance</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> Finally, we acknowledged the complexity of societies and the need for an interdisciplinary approach to comprehending their evolution and challenges</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> Through this exploration, we gained valuable insights into the foundations of human civilization and the forces that drive social and political change throughout history</w:t>
      </w:r>
      <w:r w:rsidR="005838B3">
        <w:rPr>
          <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p/><w:p/>
'@

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
